# Scheduled market data refresh: update currentAveragePrice / Leve price / profit columns
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 7758.3125
$ws.Range("I86").Value = 1628.1111
$ws.Range("K86").Value = 1628.1111
$ws.Range("M86").Value = -505.1111000000001

# Row 89
$ws.Range("H89").Value = 7758.3125
$ws.Range("I89").Value = 1628.1111
$ws.Range("K89").Value = 8140.5555
$ws.Range("M89").Value = -2524.5555

# Row 120
$ws.Range("H120").Value = 45761
$ws.Range("J120").Value = 45761
$ws.Range("L120").Value = 45761
$ws.Range("N120").Value = -55437

# Row 129
$ws.Range("H129").Value = 854.7222
$ws.Range("J129").Value = 857.94116
$ws.Range("L129").Value = 2573.82348
$ws.Range("N129").Value = -12573.82348

# Row 132
$ws.Range("H132").Value = 3180.1785
$ws.Range("I132").Value = 3381.6365
$ws.Range("J132").Value = 2441.5
$ws.Range("K132").Value = 10144.9095
$ws.Range("L132").Value = 7324.5
$ws.Range("M132").Value = -7614.9095
$ws.Range("N132").Value = -12384.5

# Row 135
$ws.Range("H135").Value = 22735392
$ws.Range("I135").Value = 731.82355
$ws.Range("K135").Value = 6586.41195
$ws.Range("M135").Value = -4051.41195

# Row 138
$ws.Range("H138").Value = 1870.2909
$ws.Range("I138").Value = 556.76
$ws.Range("J138").Value = 2964.9
$ws.Range("K138").Value = 1670.28
$ws.Range("L138").Value = 8894.700000000001
$ws.Range("M138").Value = 3469.72
$ws.Range("N138").Value = -19174.7

# Row 141
$ws.Range("H141").Value = 3011.875
$ws.Range("I141").Value = 2019
$ws.Range("K141").Value = 6057
$ws.Range("M141").Value = -877

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 49578.727
$ws.Range("I32").Value = 56096.42
$ws.Range("J32").Value = 8300
$ws.Range("K32").Value = 56096.42
$ws.Range("L32").Value = 8300
$ws.Range("M32").Value = -55809.42
$ws.Range("N32").Value = -8874

# Row 61
$ws.Range("H61").Value = 2300.2173
$ws.Range("I61").Value = 1339.1666
$ws.Range("J61").Value = 5760
$ws.Range("K61").Value = 1339.1666
$ws.Range("L61").Value = 5760
$ws.Range("M61").Value = -1127.1666
$ws.Range("N61").Value = -6184

# Row 97
$ws.Range("H97").Value = 2394.5833
$ws.Range("I97").Value = 2006.1111
$ws.Range("J97").Value = 3560
$ws.Range("K97").Value = 2006.1111
$ws.Range("L97").Value = 3560
$ws.Range("M97").Value = -1510.1111
$ws.Range("N97").Value = -4552

# Row 136
$ws.Range("H136").Value = 2300.2173
$ws.Range("I136").Value = 1339.1666
$ws.Range("J136").Value = 5760
$ws.Range("K136").Value = 4017.4998
$ws.Range("L136").Value = 17280
$ws.Range("M136").Value = -1467.4998
$ws.Range("N136").Value = -22380

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1776.1111
$ws.Range("I20").Value = 2094.1667
$ws.Range("J20").Value = 1140
$ws.Range("K20").Value = 2094.1667
$ws.Range("L20").Value = 1140
$ws.Range("M20").Value = -1847.1667
$ws.Range("N20").Value = -1634

# Row 94
$ws.Range("H94").Value = 3103.0322
$ws.Range("I94").Value = 1714.5
$ws.Range("J94").Value = 5627.636
$ws.Range("K94").Value = 1714.5
$ws.Range("L94").Value = 5627.636
$ws.Range("M94").Value = -1263.5
$ws.Range("N94").Value = -6529.636

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2560.9666
$ws.Range("I31").Value = 2182.9
$ws.Range("J31").Value = 3317.1
$ws.Range("K31").Value = 2182.9
$ws.Range("L31").Value = 3317.1
$ws.Range("M31").Value = -1887.9
$ws.Range("N31").Value = -3907.1

# Row 34
$ws.Range("H34").Value = 2560.9666
$ws.Range("I34").Value = 2182.9
$ws.Range("J34").Value = 3317.1
$ws.Range("K34").Value = 2182.9
$ws.Range("L34").Value = 3317.1
$ws.Range("M34").Value = -1980.9
$ws.Range("N34").Value = -3721.1

# Row 62
$ws.Range("H62").Value = 6002.8335
$ws.Range("I62").Value = 4002.5
$ws.Range("K62").Value = 4002.5
$ws.Range("M62").Value = -3378.5

# Row 65
$ws.Range("H65").Value = 6002.8335
$ws.Range("I65").Value = 4002.5
$ws.Range("K65").Value = 20012.5
$ws.Range("M65").Value = -16892.5

# Row 86
$ws.Range("H86").Value = 15518.909
$ws.Range("I86").Value = 2610
$ws.Range("J86").Value = 26276.334
$ws.Range("K86").Value = 2610
$ws.Range("L86").Value = 26276.334
$ws.Range("M86").Value = -1487
$ws.Range("N86").Value = -28522.334

# Row 89
$ws.Range("H89").Value = 15518.909
$ws.Range("I89").Value = 2610
$ws.Range("J89").Value = 26276.334
$ws.Range("K89").Value = 13050
$ws.Range("L89").Value = 131381.67
$ws.Range("M89").Value = -7434
$ws.Range("N89").Value = -142613.67

# Row 99
$ws.Range("H99").Value = 14586521
$ws.Range("I99").Value = 2978594.2
$ws.Range("J99").Value = 41671684
$ws.Range("K99").Value = 2978594.2
$ws.Range("L99").Value = 41671684
$ws.Range("M99").Value = -2977096.2
$ws.Range("N99").Value = -41674680

# Row 126
$ws.Range("H126").Value = 14586521
$ws.Range("I126").Value = 2978594.2
$ws.Range("J126").Value = 41671684
$ws.Range("K126").Value = 8935782.600000001
$ws.Range("L126").Value = 125015052
$ws.Range("M126").Value = -8933312.600000001
$ws.Range("N126").Value = -125019992

# Row 134
$ws.Range("H134").Value = 1101.3182
$ws.Range("I134").Value = 806.9286
$ws.Range("J134").Value = 1616.5
$ws.Range("K134").Value = 2420.7858
$ws.Range("L134").Value = 4849.5
$ws.Range("M134").Value = 114.2142000000003
$ws.Range("N134").Value = -9919.5

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 265.8
$ws.Range("I33").Value = 199
$ws.Range("J33").Value = 282.5
$ws.Range("K33").Value = 1194
$ws.Range("L33").Value = 1695
$ws.Range("M33").Value = -911
$ws.Range("N33").Value = -2261

# Row 117
$ws.Range("H117").Value = 2142.25
$ws.Range("I117").Value = 741.2
$ws.Range("J117").Value = 4477.3335
$ws.Range("K117").Value = 2223.6
$ws.Range("L117").Value = 13432.0005
$ws.Range("M117").Value = 1218.4
$ws.Range("N117").Value = -20316.0005

# Row 131
$ws.Range("H131").Value = 757.4
$ws.Range("I131").Value = 354
$ws.Range("J131").Value = 778.6316
$ws.Range("K131").Value = 1062
$ws.Range("L131").Value = 2335.8948
$ws.Range("M131").Value = 3978
$ws.Range("N131").Value = -12415.8948

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 87.083336
$ws.Range("I2").Value = 47.285713
$ws.Range("J2").Value = 142.8
$ws.Range("K2").Value = 47.285713
$ws.Range("L2").Value = 142.8
$ws.Range("M2").Value = 65.714287
$ws.Range("N2").Value = -368.8

# Row 97
$ws.Range("H97").Value = 3781.375
$ws.Range("I97").Value = 1729.1666
$ws.Range("K97").Value = 1729.1666
$ws.Range("M97").Value = -1233.1666

# Row 102
$ws.Range("H102").Value = 1554.1428
$ws.Range("I102").Value = 1603.9656
$ws.Range("J102").Value = 1313.3334
$ws.Range("K102").Value = 1603.9656
$ws.Range("L102").Value = 1313.3334
$ws.Range("M102").Value = 18.03440000000001
$ws.Range("N102").Value = -4557.3334

# Row 113
$ws.Range("H113").Value = 3340.25
$ws.Range("J113").Value = 3944.4
$ws.Range("L113").Value = 3944.4
$ws.Range("N113").Value = -8284.4

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 3707.875
$ws.Range("J61").Value = 6988.3335
$ws.Range("L61").Value = 6988.3335
$ws.Range("N61").Value = -7392.3335

# Row 68
$ws.Range("H68").Value = 5138.1816
$ws.Range("I68").Value = 2306
$ws.Range("J68").Value = 7498.3335
$ws.Range("K68").Value = 2306
$ws.Range("L68").Value = 7498.3335
$ws.Range("M68").Value = -1557
$ws.Range("N68").Value = -8996.333500000001

# Row 71
$ws.Range("H71").Value = 5138.1816
$ws.Range("I71").Value = 2306
$ws.Range("J71").Value = 7498.3335
$ws.Range("K71").Value = 11530
$ws.Range("L71").Value = 37491.6675
$ws.Range("M71").Value = -7786
$ws.Range("N71").Value = -44979.6675

# Row 104
$ws.Range("H104").Value = 23937.555
$ws.Range("J104").Value = 23937.555
$ws.Range("L104").Value = 23937.555
$ws.Range("N104").Value = -30925.555

# Row 113
$ws.Range("H113").Value = 3707.875
$ws.Range("J113").Value = 6988.3335
$ws.Range("L113").Value = 6988.3335
$ws.Range("N113").Value = -11328.3335

$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 20249
$ws.Range("J63").Value = 20249
$ws.Range("L63").Value = 20249
$ws.Range("N63").Value = -21497

# Row 66
$ws.Range("H66").Value = 20249
$ws.Range("J66").Value = 20249
$ws.Range("L66").Value = 60747
$ws.Range("N66").Value = -66987

# Row 132
$ws.Range("H132").Value = 2075.25
$ws.Range("I132").Value = 920.8
$ws.Range("K132").Value = 2762.4
$ws.Range("M132").Value = -232.3999999999996

# Row 136
$ws.Range("H136").Value = 18519764
$ws.Range("I136").Value = 26316886
$ws.Range("J136").Value = 1603.625
$ws.Range("K136").Value = 78950658
$ws.Range("L136").Value = 4810.875
$ws.Range("M136").Value = -78948108
$ws.Range("N136").Value = -9910.875
